$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.108.86'
$ws.Range("E2").Value = '  -2.46%  '

$ws.Range("D3").Value = '1.898.15'
$ws.Range("E3").Value = '  -2.83%  '

$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '332.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.06%  '

$ws.Range("E6").Value = '  +0.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4594'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.83%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4131'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.33%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.85'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.63%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07997'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.36%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.011'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.81%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.18'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.54%  '

$ws.Range("D13").Value = '1.886.56'
$ws.Range("E13").Value = '  -3.27%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.939'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.20%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.127'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.18%  '

$ws.Range("E16").Value = '  +0.18%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.05'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.60%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001029'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.15%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06565'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.89%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.62'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.70%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.03%  '

$ws.Range("D22").Value = '29.043.46'
$ws.Range("E22").Value = '  -2.54%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.497'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.77%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.43'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.203'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.48%  '

$ws.Range("D26").Value = '2.122.13'
$ws.Range("E26").Value = '  -2.67%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.96%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.75'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.42%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.121'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.20%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.637'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.94%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.11'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.87%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.048'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09381'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.75%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.418'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.26%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.534'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.14%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.354'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.37%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06096'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.70%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02238'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.07%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.433'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.32%  '

$ws.Range("E40").Value = '  -1.53%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5834'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.45%  '

$ws.Range("E42").Value = '  +0.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1829'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.56%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.13'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.65%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.248'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.81%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07521'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.23%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.303'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.46%  '

$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5526'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.62%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '12.06'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.58%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.924'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.48%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '112.21'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.77%  '
